$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.579.71'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '2.219.21'
$ws.Range("E3").Value = '  -1.02%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").Value = '250.11'
$ws.Range("E5").Value = '  +7.72%  '

$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("D7").Value = '70.38'
$ws.Range("E7").Value = '  +1.74%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  +3.69%  '

$ws.Range("D10").Value = '41.50'
$ws.Range("E10").Value = '  +17.01%  '

$ws.Range("D11").Value = '0.0958'
$ws.Range("E11").Value = '  -2.44%  '

$ws.Range("D12").Value = '58.54'
$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("E14").Value = '  +3.50%  '

$ws.Range("D15").Value = '2.550.01'
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("D16").Value = '14.83'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").Value = '0.852'
$ws.Range("E17").Value = '  -0.17%  '

$ws.Range("D18").Value = '2.219.32'
$ws.Range("E18").Value = '  -0.90%  '

$ws.Range("D19").Value = '41.469.36'
$ws.Range("E19").Value = '  -0.89%  '

$ws.Range("D20").Value = '0.0₃0964'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  -0.61%  '

$ws.Range("D22").Value = '72.44'
$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").Value = '233.68'
$ws.Range("E23").Value = '  -0.87%  '

$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +8.57%  '

$ws.Range("D25").Value = '3.88'
$ws.Range("E25").Value = '  +6.68%  '

$ws.Range("D27").Value = '2.50'
$ws.Range("E27").Value = '  +6.67%  '

$ws.Range("D28").Value = '10.48'
$ws.Range("E28").Value = '  +5.18%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  +0.59%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '171.13'
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("D31").Value = '20.56'
$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").Value = '0.120'
$ws.Range("E32").Value = '  +1.50%  '

$ws.Range("E33").Value = '  -2.18%  '

$ws.Range("E34").Value = '  +1.61%  '

$ws.Range("D35").Value = '0.0716'
$ws.Range("E35").Value = '  +0.67%  '

$ws.Range("D36").Value = '4.66'
$ws.Range("E36").Value = '  -1.67%  '

$ws.Range("D37").Value = '25.94'
$ws.Range("E37").Value = '  +17.86%  '

$ws.Range("D38").Value = '3.94'
$ws.Range("E38").Value = '  +10.17%  '

$ws.Range("D39").Value = '0.0291'
$ws.Range("E39").Value = '  +10.21%  '

$ws.Range("D40").Value = '2.28'
$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("D41").Value = '68.70'
$ws.Range("E41").Value = '  +3.55%  '

$ws.Range("D42").Value = '5.91'
$ws.Range("E42").Value = '  -1.41%  '

$ws.Range("D43").Value = '11.81'
$ws.Range("E43").Value = '  +19.53%  '

$ws.Range("D44").Value = '0.207'
$ws.Range("E44").Value = '  +9.35%  '

$ws.Range("D45").Value = '4.90'
$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("B46").Value = 'SynthetixNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D46").Value = '4.77'
$ws.Range("E46").Value = '  +11.88%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '8.70'
$ws.Range("E47").Value = '  -3.53%  '

$ws.Range("E48").Value = '  +0.80%  '

$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("E50").Value = '  +6.81%  '

$ws.Range("D51").Value = '1.19'
$ws.Range("E51").Value = '  +1.65%  '
